$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Update existing rows 33-35: Results column (E) PASS -> SKIP
$ws.Range("E33").Value = "SKIP"
$ws.Range("E34").Value = "SKIP"
$ws.Range("E35").Value = "SKIP"

# Copy formatting of row 35 down into the two new rows (36-37) so the new
# cells pick up the same cell styles used by the rest of the table.
$ws.Range("A35:E35").Copy()
$ws.Range("A36:E36").PasteSpecial(-4122)
$ws.Range("A35:E35").Copy()
$ws.Range("A37:E37").PasteSpecial(-4122)

# Row 36 - new test case: PublishedAPostLikeCountTest
$ws.Range("C36").Value = "Verify that POST tab count getting increased while appreciate post from Record view page"
$ws.Range("A36").Value = "PublishedAPostLikeCountTest"
$ws.Range("B36").Value = "TBD"
$ws.Range("D36").Value = "Y"
$ws.Range("E36").Value = "SKIP"

# Row 37 - new test case: PublishedAPostTimeStampTest
$ws.Range("A37").Value = "PublishedAPostTimeStampTest"
$ws.Range("C37").Value = "Verify that Created Post displayed as per System date"
$ws.Range("B37").Value = "TBD"
$ws.Range("D37").Value = "Y"
$ws.Range("E37").Value = "PASS"

# Refresh the view to show the newly added rows
$ws.Range("D2:D37").Select()
$excel.ActiveWindow.ScrollRow = 12
$excel.ActiveWindow.ScrollColumn = 1
